$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("Assignment" and everything to its
# right shifts one column to the right: D->E, E->F, F->G).
$ws.Columns("D").Insert()

# Set the new column D width (~18.33 characters).
$ws.Columns("D").ColumnWidth = 17.5

# Header for the new "Slides" column.
$ws.Range("D1").Value = "Slides"

# Populate the new Slides column for the first two sessions.
$ws.Range("D2").Value = "01_introduction"
$ws.Range("D3").Value = "02_data"

# Update the content description for the second session.
$ws.Range("C3").Value = "Data sources and definitions"

# Populate the Slides column for the third session.
$ws.Range("D4").Value = "03_evolution"

# Update the active cell selection.
$ws.Range("D5").Select()
